$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from its old location (between the
#    "...keterangan_cerai} " run and the "${lanjutan}" run) to the
#    empty paragraph that immediately follows "MENUGASKAN :".
# ------------------------------------------------------------------

# Remove the existing _GoBack bookmark (it sits mid-paragraph further
# down in the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate "MENUGASKAN :" and grab the paragraph right after it -- using
# Find + relative paragraph navigation (.Next()) rather than an
# absolute Paragraphs(n) index, which keeps the anchor accurate even
# though the body also contains tables.
$findRange = $d.Content
$findRange.Find.Execute("MENUGASKAN :") | Out-Null
$findRange.Collapse(0)
$menugaskanPara = $findRange.Paragraphs(1)
$targetPara = $menugaskanPara.Next()

$d.Bookmarks.Add("_GoBack", $targetPara.Range) | Out-Null

# ------------------------------------------------------------------
# 2) Drop the duplicate empty paragraph (the one with
#    ind left="993" firstLine="992") that sits right before the
#    signature table, keeping the first empty paragraph intact.
# ------------------------------------------------------------------

$findRange2 = $d.Content
$findRange2.Find.Execute("Kepada pihak-pihak yang bersangkutan") | Out-Null
$findRange2.Collapse(0)
$closingPara = $findRange2.Paragraphs(1)

$tabsPara = $closingPara.Next()     # paragraph full of tab characters
$paraKeep = $tabsPara.Next()        # first empty paragraph (kept)
$paraDrop = $paraKeep.Next()        # second empty paragraph (removed)

$paraDrop.Range.Delete() | Out-Null

# ------------------------------------------------------------------
# 3) Fix the first cell width of the signature table: 5399 -> 5400 dxa
#    (dxa / 20 = points).
# ------------------------------------------------------------------

$signatureTable = $d.Tables(2)
$signatureTable.Cell(1, 1).Width = 5400 / 20
